# Fix trajectory read for real robot data:
#  - Column B (col 2): sign flipped (negate)
#  - Column C (col 3): replaced by (pi/2 - value)
#  - Column D (col 4): sign flipped (negate)
#  - Column E (col 5): replaced by (pi/2 - value)
#  - Column F (col 6): sign flipped (negate)
# Columns A and G are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$halfPi = 1.5707963267948966

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $bCell.Value2 = -1 * $bCell.Value2
    $cCell.Value2 = $halfPi - $cCell.Value2
    $dCell.Value2 = -1 * $dCell.Value2
    $eCell.Value2 = $halfPi - $eCell.Value2
    $fCell.Value2 = -1 * $fCell.Value2
}
